$d = $word.ActiveDocument
$p = $d.Paragraphs
$p.Item(1).Range.Text = "2023-06-12 Monday"
$p.Item(2).Range.Text = "76-17=59"
$p.Item(3).Range.Text = "52-3=49"
$p.Item(4).Range.Text = "92-34=58"
$p.Item(5).Range.Text = "12+43=55"
$p.Item(6).Range.Text = "22+12=34"
$p.Item(8).Range.Text = "24-8=16"
$p.Item(9).Range.Text = "35+13=48"
$p.Item(10).Range.Text = "40-6=34"
$p.Item(11).Range.Text = "70-44=26"
$p.Item(12).Range.Text = "48+23=71"
$p.Item(14).Range.Text = "32+66=98"
$p.Item(15).Range.Text = "45+11=56"
$p.Item(16).Range.Text = "45-19=26"
$p.Item(17).Range.Text = "58-15=43"
$p.Item(18).Range.Text = "3+65=68"
$p.Item(20).Range.Text = "11+15=26"
$p.Item(21).Range.Text = "28+51=79"
$p.Item(22).Range.Text = "18+22=40"
$p.Item(23).Range.Text = "4+5=9"
$p.Item(24).Range.Text = "67+26=93"
$p.Item(26).Range.Text = "79-46=33"
$p.Item(27).Range.Text = "55-40=15"
$p.Item(28).Range.Text = "96-90=6"
$p.Item(29).Range.Text = "26+52=78"
$p.Item(30).Range.Text = "8+29=37"
$p.Item(32).Range.Text = "10+22=32"
$p.Item(33).Range.Text = "34-31=3"
$p.Item(34).Range.Text = "43+1=44"
$p.Item(35).Range.Text = "54+40=94"
$p.Item(36).Range.Text = "99-25=74"
$p.Item(38).Range.Text = "97-72=25"
$p.Item(39).Range.Text = "97-3=94"
$p.Item(40).Range.Text = "33+16=49"
$p.Item(41).Range.Text = "95-33=62"
$p.Item(42).Range.Text = "97+2=99"
$p.Item(44).Range.Text = "65-26=39"
$p.Item(45).Range.Text = "74-33=41"
$p.Item(46).Range.Text = "69-52=17"
$p.Item(47).Range.Text = "9+74=83"
$p.Item(48).Range.Text = "91-89=2"
$p.Item(50).Range.Text = "84-1=83"
$p.Item(51).Range.Text = "49+37=86"
$p.Item(52).Range.Text = "23-23=0"
$p.Item(53).Range.Text = "61-57=4"
$p.Item(54).Range.Text = "88-88=0"
$p.Item(56).Range.Text = "81-10=71"
$p.Item(57).Range.Text = "83-82=1"
$p.Item(58).Range.Text = "18+9=27"
$p.Item(59).Range.Text = "47-7=40"
$p.Item(60).Range.Text = "52+28=80"
$p.Item(62).Range.Text = "17+49=66"
$p.Item(63).Range.Text = "94-57=37"
$p.Item(64).Range.Text = "25+6=31"
$p.Item(65).Range.Text = "74+23=97"
$p.Item(66).Range.Text = "48+48=96"
$p.Item(68).Range.Text = "14+9=23"
$p.Item(69).Range.Text = "98-65=33"
$p.Item(70).Range.Text = "50+19=69"
$p.Item(71).Range.Text = "29+35=64"
$p.Item(72).Range.Text = "90-69=21"
$p.Item(74).Range.Text = "84-21=63"
$p.Item(75).Range.Text = "5+26=31"
$p.Item(76).Range.Text = "74-20=54"
$p.Item(77).Range.Text = "29-28=1"
$p.Item(78).Range.Text = "76-28=48"
$p.Item(80).Range.Text = "53+15=68"
$p.Item(81).Range.Text = "7+62=69"
$p.Item(82).Range.Text = "93-56=37"
$p.Item(83).Range.Text = "78-70=8"
$p.Item(84).Range.Text = "98-43=55"
$p.Item(86).Range.Text = "98-12=86"
$p.Item(87).Range.Text = "37+31=68"
$p.Item(88).Range.Text = "96-54=42"
$p.Item(89).Range.Text = "1+86=87"
$p.Item(90).Range.Text = "98-46=52"
$p.Item(92).Range.Text = "26-0=26"
$p.Item(93).Range.Text = "90-74=16"
$p.Item(94).Range.Text = "18+15=33"
$p.Item(95).Range.Text = "37+38=75"
$p.Item(96).Range.Text = "27+7=34"
$p.Item(98).Range.Text = "45+41=86"
$p.Item(99).Range.Text = "58+33=91"
$p.Item(100).Range.Text = "14+37=51"
$p.Item(101).Range.Text = "29+23=52"
$p.Item(102).Range.Text = "69+7=76"
$p.Item(104).Range.Text = "70+8=78"
$p.Item(105).Range.Text = "62-33=29"
$p.Item(106).Range.Text = "91-15=76"
$p.Item(107).Range.Text = "54-11=43"
$p.Item(108).Range.Text = "24-12=12"
$p.Item(110).Range.Text = "61-31=30"
$p.Item(111).Range.Text = "31+35=66"
$p.Item(112).Range.Text = "42+36=78"
$p.Item(113).Range.Text = "75-53=22"
$p.Item(114).Range.Text = "69-18=51"
$p.Item(116).Range.Text = "35+23=58"
$p.Item(117).Range.Text = "75+18=93"
$p.Item(118).Range.Text = "62-17=45"
$p.Item(119).Range.Text = "18+15=33"
$p.Item(120).Range.Text = "22-5=17"
